$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to English snake_case names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the Spanish connector words (de/del/la/las/el/los) in
# state/municipality names so they read 'De'/'Del'/'La'/... instead of lowercase
$ws.Range("B4").Value = "Amatenango De La Frontera"
$ws.Range("B10").Value = "Hidalgo Del Parral"
$ws.Range("A15").Value = "Ciudad De México"
$ws.Range("A23").Value = "Estado De México"
$ws.Range("B23").Value = "Atizapán De Zaragoza"
$ws.Range("B25").Value = "Ecatepec De Morelos"
$ws.Range("B30").Value = "San Martín De Las Pirámides"
$ws.Range("B33").Value = "Tlalnepantla De Baz"
$ws.Range("B36").Value = "Apaseo El Alto"
$ws.Range("B38").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B42").Value = "Acapulco De Juárez"
$ws.Range("B43").Value = "Alcozauca De Guerrero"
$ws.Range("B45").Value = "Buenavista De Cuéllar"
$ws.Range("B46").Value = "Chilpancingo De Los Bravo"
$ws.Range("B48").Value = "Iguala De La Independencia"
$ws.Range("B49").Value = "Zihuatanejo De Azueta"
$ws.Range("B54").Value = "Tepecoacuilco De Trujano"
$ws.Range("B61").Value = "Huasca De Ocampo"
$ws.Range("B63").Value = "Mineral De La Reforma"
$ws.Range("B64").Value = "Mixquiahuala De Juárez"
$ws.Range("B65").Value = "Nopala De Villagrán"
$ws.Range("B66").Value = "Omitlán De Juárez"
$ws.Range("B67").Value = "Pachuca De Soto"
$ws.Range("B72").Value = "Tula De Allende"
$ws.Range("B109").Value = "Izúcar De Matamoros"
$ws.Range("B118").Value = "Landa De Matamoros"
$ws.Range("B142").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B145").Value = "Ignacio De La Llave"
$ws.Range("B150").Value = "Ozuluama De Mascareñas"
$ws.Range("B151").Value = "Poza Rica De Hidalgo"

# Drop the stray footnote/source rows below the data table (rows 165-169 and
# 476-480) - clearing them removes the rows from the sheet and the used range
# (dimension) shrinks back down to A1:D163 automatically.
$ws.Rows.Item(165).ClearContents()
$ws.Rows.Item(166).ClearContents()
$ws.Rows.Item(167).ClearContents()
$ws.Rows.Item(168).ClearContents()
$ws.Rows.Item(169).ClearContents()
$ws.Rows.Item(476).ClearContents()
$ws.Rows.Item(477).ClearContents()
$ws.Rows.Item(478).ClearContents()
$ws.Rows.Item(479).ClearContents()
$ws.Rows.Item(480).ClearContents()
